$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.836.56'
$ws.Range('E2').Value = '  +7.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.696.56'
$ws.Range('E3').Value = '  +11.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '515.18'
$ws.Range('E5').Value = '  +6.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.87'
$ws.Range('E6').Value = '  +3.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.994'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.691.67'
$ws.Range('E9').Value = '  +10.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.20'
$ws.Range('E10').Value = '  +10.13%  '
$ws.Range('E11').Value = '  +6.68%  '
$ws.Range('E12').Value = '  +5.19%  '
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.124.54'
$ws.Range('E14').Value = '  +10.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '61.019.91'
$ws.Range('E15').Value = '  +7.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.12'
$ws.Range('E16').Value = '  +7.06%  '
$ws.Range('E17').Value = '  +6.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.678.74'
$ws.Range('E18').Value = '  +10.10%  '
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '350.20'
$ws.Range('E20').Value = '  +7.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.57'
$ws.Range('E21').Value = '  +6.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.23'
$ws.Range('E22').Value = '  +5.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.38'
$ws.Range('E24').Value = '  +4.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.426'
$ws.Range('E25').Value = '  +4.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.764.96'
$ws.Range('E26').Value = '  +9.32%  '
$ws.Range('E27').Value = '  +6.16%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0875'
$ws.Range('E29').Value = '  +12.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.60'
$ws.Range('E30').Value = '  +5.26%  '
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.67'
$ws.Range('E32').Value = '  +5.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.17'
$ws.Range('E33').Value = '  +5.64%  '
$ws.Range('E34').Value = '  +4.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.75'
$ws.Range('E35').Value = '  +8.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.08'
$ws.Range('E36').Value = '  +10.77%  '
$ws.Range('E37').Value = '  +7.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.890'
$ws.Range('E38').Value = '  +5.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('E39').Value = '  +13.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '311.19'
$ws.Range('E40').Value = '  +17.11%  '
$ws.Range('E41').Value = '  +8.49%  '
$ws.Range('E42').Value = '  +31.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.72'
$ws.Range('E43').Value = '  +4.66%  '
$ws.Range('E44').Value = '  +9.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0580'
$ws.Range('E45').Value = '  +9.56%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.46'
$ws.Range('E47').Value = '  +17.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.06'
$ws.Range('E48').Value = '  +8.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.996'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0238'
$ws.Range('E50').Value = '  +4.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.049.01'
$ws.Range('E51').Value = '  +10.44%  '
